$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$changes = @{
    "H4" = 449.66666
    "I4" = 149.33333
    "J4" = 750.0
    "K4" = 149.33333
    "L4" = 750.0
    "M4" = -35.33332999999999
    "N4" = -978.0
    "H12" = 286.4
    "I12" = 235.5
    "K12" = 235.5
    "M12" = -65.5
    "H17" = 2364.0
    "J17" = 2364.0
    "L17" = 7092.0
    "N17" = -7428.0
    "H33" = 3250.9707
    "I33" = 4024.4614
    "K33" = 4024.4614
    "M33" = -3795.4614
    "H40" = 4288.0
    "J40" = 4872.5
    "L40" = 4872.5
    "N40" = -5222.5
    "H48" = 3035.5715
    "J48" = 3333.1667
    "L48" = 9999.500100000001
    "N48" = -10583.5001
    "H51" = 9799.444
    "J51" = 8099.0
    "L51" = 8099.0
    "N51" = -9067.0
    "H56" = 3035.5715
    "J56" = 3333.1667
    "L56" = 9999.500100000001
    "N56" = -11067.5001
    "H64" = 3449.0
    "I64" = 3449.0
    "J64" = 0.0
    "K64" = 3449.0
    "L64" = 0.0
    "M64" = $null
    "N64" = -3201.0
    "H67" = 3449.0
    "I67" = 3449.0
    "J67" = 0.0
    "K67" = 3449.0
    "L67" = 0.0
    "M67" = $null
    "N67" = -2591.0
    "H92" = 100657.7
    "I92" = 125172.125
    "J92" = 2600.0
    "K92" = 125172.125
    "L92" = 2600.0
    "M92" = -123924.125
    "N92" = -5096.0
    "H98" = 2340.0557
    "I98" = 2340.0557
    "K98" = 2340.0557
    "M98" = -842.0556999999999
    "H106" = 5240.1763
    "I106" = 3367.2
    "K106" = 3367.2
    "M106" = -2736.2
    "H107" = 2012.3
    "I107" = 1891.125
    "K107" = 1891.125
    "M107" = 28.875
    "H112" = 3430.923
    "J112" = 3430.923
    "L112" = 10292.769
    "N112" = -12508.769
    "H122" = 2340.0557
    "I122" = 2340.0557
    "K122" = 7020.1671
    "M122" = -4570.1671
    "H129" = 984.26666
    "I129" = 943.38464
    "K129" = 2830.15392
    "M129" = 2169.84608
    "H132" = 2274.2856
    "I132" = 1820.0883
    "K132" = 5460.2649
    "M132" = -2930.2649
    "H137" = 42141.16
    "I137" = 73063.57
    "K137" = 219190.71
    "M137" = -216640.71
    "H138" = 2001.9375
    "I138" = 1975.4333
    "J138" = 2399.5
    "K138" = 5926.2999
    "L138" = 7198.5
    "M138" = -786.2999
    "N138" = -17478.5
    "H141" = 1188.3125
    "I141" = 1251.9333
    "J141" = 234.0
    "K141" = 3755.7999
    "L141" = 702.0
    "M141" = 1424.2001
    "N141" = -11062.0
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}

$ws = $wb.Worksheets.Item("ARM")
$changes = @{
    "H2" = 2653.6924
    "I2" = 1944.4445
    "K2" = 1944.4445
    "M2" = -1831.4445
    "H74" = 25045.342
    "I74" = 27956.691
    "J74" = 2336.8
    "K74" = 27956.691
    "L74" = 2336.8
    "M74" = -27082.691
    "N74" = -4084.8
    "H77" = 25045.342
    "I77" = 27956.691
    "J77" = 2336.8
    "K77" = 139783.455
    "L77" = 11684.0
    "M77" = -135415.455
    "N77" = -20420.0
    "H102" = 2685.2
    "I102" = 2454.8333
    "K102" = 2454.8333
    "M102" = -832.8332999999998
    "H116" = 2653.6924
    "I116" = 1944.4445
    "K116" = 1944.4445
    "M116" = 349.5554999999999
    "H122" = 2512.5625
    "I122" = 2442.9285
    "J122" = 3000.0
    "K122" = 7328.7855
    "L122" = 9000.0
    "M122" = -4878.7855
    "N122" = -13900.0
    "H132" = 32983.484
    "I132" = 32983.484
    "J132" = 0.0
    "K132" = 98950.45199999999
    "L132" = 0.0
    "M132" = $null
    "N132" = -96420.45199999999
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}

$ws = $wb.Worksheets.Item("BSM")
$changes = @{
    "H3" = 2653.6924
    "I3" = 1944.4445
    "K3" = 1944.4445
    "M3" = -1830.4445
    "H105" = 2693.4375
    "I105" = 2215.2307
    "K105" = 2215.2307
    "M105" = -468.2307000000001
    "H134" = 1664.72
    "I134" = 1664.72
    "J134" = 0.0
    "K134" = 4994.16
    "L134" = 0.0
    "M134" = $null
    "N134" = -2459.16
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CRP")
$changes = @{
    "H31" = 3714.1365
    "I31" = 3095.1667
    "K31" = 3095.1667
    "M31" = -2800.1667
    "H34" = 3714.1365
    "I34" = 3095.1667
    "K34" = 3095.1667
    "M34" = -2893.1667
    "H86" = 3401.375
    "I86" = 3882.6667
    "J86" = 3112.6
    "K86" = 3882.6667
    "L86" = 3112.6
    "M86" = -2759.6667
    "N86" = -5358.6
    "H89" = 3401.375
    "I89" = 3882.6667
    "J89" = 3112.6
    "K89" = 19413.3335
    "L89" = 15563.0
    "M89" = -13797.3335
    "N89" = -26795.0
    "H105" = 50986.75
    "I105" = 50986.75
    "K105" = 50986.75
    "M105" = -49239.75
    "H121" = 49663.0
    "J121" = 49663.0
    "L121" = 49663.0
    "N121" = -52283.0
    "H122" = 1528.7693
    "I122" = 1489.5834
    "K122" = 4468.7502
    "M122" = -2018.7502
    "H125" = 46316.5
    "J125" = 46316.5
    "L125" = 46316.5
    "N125" = -51236.5
    "H132" = 977.7027
    "I132" = 1013.46875
    "J132" = 748.8
    "K132" = 3040.40625
    "L132" = 2246.4
    "M132" = -510.40625
    "N132" = -7306.4
    "H134" = 28959.555
    "I134" = 28959.555
    "K134" = 86878.66500000001
    "M134" = -84343.66500000001
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CUL")
$changes = @{
    "H2" = 55.64706
    "J2" = 62.0
    "L2" = 372.0
    "N2" = -598.0
    "H5" = 961.53845
    "J5" = 1335.6
    "L5" = 4006.8
    "N5" = -4230.799999999999
    "H12" = 215.58824
    "I12" = 195.6
    "J12" = 223.91667
    "K12" = 586.8
    "L12" = 671.75001
    "M12" = -413.8
    "N12" = -1017.75001
    "H37" = 64243.75
    "J37" = 64243.75
    "L37" = 192731.25
    "N37" = -192955.25
    "H80" = 123079.0
    "I80" = 299999.0
    "J80" = 78849.0
    "K80" = 899997.0
    "L80" = 236547.0
    "M80" = -899061.0
    "N80" = -238419.0
    "H83" = 123079.0
    "I83" = 299999.0
    "J83" = 78849.0
    "K83" = 2699991.0
    "L83" = 709641.0
    "M83" = -2695311.0
    "N83" = -719001.0
    "H135" = 961.53845
    "J135" = 1335.6
    "L135" = 12020.4
    "N135" = -17090.4
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}

$ws = $wb.Worksheets.Item("GSM")
$changes = @{
    "H3" = 300.0
    "I3" = 300.0
    "K3" = 300.0
    "M3" = -184.0
    "H15" = 50000.0
    "J15" = 50000.0
    "L15" = 50000.0
    "N15" = -50576.0
    "H20" = 0.0
    "J20" = 0.0
    "L20" = $null
    "N20" = 0.0
    "H81" = 50000.0
    "J81" = 50000.0
    "L81" = 50000.0
    "N81" = -51996.0
    "H84" = 50000.0
    "J84" = 50000.0
    "L84" = 150000.0
    "N84" = -159984.0
    "H107" = 31361.697
    "I107" = 53398.42
    "J107" = 1454.7142
    "K107" = 53398.42
    "L107" = 1454.7142
    "M107" = -51478.42
    "N107" = -5294.7142
    "H122" = 7412.4614
    "I122" = 6640.2173
    "J122" = 13333.0
    "K122" = 19920.6519
    "L122" = 39999.0
    "M122" = -17470.6519
    "N122" = -44899.0
    "H132" = 32397.158
    "I132" = 39691.934
    "J132" = 5041.75
    "K132" = 119075.802
    "L132" = 15125.25
    "M132" = -116545.802
    "N132" = -20185.25
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}

$ws = $wb.Worksheets.Item("LTW")
$changes = @{
    "H22" = 46794.32
    "I22" = 111927.3
    "J22" = 3372.3333
    "K22" = 111927.3
    "L22" = 3372.3333
    "M22" = -111632.3
    "N22" = -3962.3333
    "H23" = 13005.0
    "I23" = 13005.0
    "K23" = 13005.0
    "M23" = -12775.0
    "H24" = 110006.0
    "I24" = 110006.0
    "K24" = 110006.0
    "M24" = -109663.0
    "H25" = 7000.0
    "I25" = 6833.3335
    "J25" = 8000.0
    "K25" = 6833.3335
    "L25" = 8000.0
    "M25" = -6603.3335
    "N25" = -8460.0
    "H27" = 46794.32
    "I27" = 111927.3
    "J27" = 3372.3333
    "K27" = 111927.3
    "L27" = 3372.3333
    "M27" = -111820.3
    "N27" = -3586.3333
    "H40" = 19507.834
    "I40" = 24942.0
    "J40" = 10968.429
    "K40" = 24942.0
    "L40" = 10968.429
    "M40" = -24806.0
    "N40" = -11240.429
    "H46" = 14479.429
    "I46" = 32319.428
    "K46" = 32319.428
    "M46" = -32131.428
    "H61" = 6689.0
    "I61" = 2500.0
    "K61" = 2500.0
    "M61" = -2298.0
    "H68" = 3955.5715
    "J68" = 3933.3333
    "L68" = 3933.3333
    "N68" = -5431.3333
    "H71" = 3955.5715
    "J71" = 3933.3333
    "L71" = 19666.6665
    "N71" = -27154.6665
    "H93" = 1490.8857
    "I93" = 1492.3572
    "J93" = 1489.9048
    "K93" = 1492.3572
    "L93" = 1489.9048
    "M93" = -244.3571999999999
    "N93" = -3985.9048
    "H113" = 6689.0
    "I113" = 2500.0
    "K113" = 2500.0
    "M113" = -330.0
    "H125" = 87999.5
    "J125" = 87999.5
    "L125" = 87999.5
    "N125" = -97839.5
    "H132" = 23382.75
    "I132" = 26407.568
    "K132" = 79222.704
    "M132" = -76692.704
    "H136" = 3094.8572
    "I136" = 2774.6
    "K136" = 8323.8
    "M136" = -5773.799999999999
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}

$ws = $wb.Worksheets.Item("WVR")
$changes = @{
    "H7" = 3501.3333
    "I7" = 252.5
    "J7" = 9999.0
    "K7" = 252.5
    "L7" = 9999.0
    "M7" = -139.5
    "N7" = -10225.0
    "H8" = 6000.0
    "J8" = 6000.0
    "L8" = 6000.0
    "N8" = -6280.0
    "H21" = 0.0
    "J21" = 0.0
    "L21" = $null
    "N21" = 0.0
    "H35" = 0.0
    "J35" = 0.0
    "L35" = $null
    "N35" = 0.0
    "J54" = 0.0
    "L54" = 0.0
    "N54" = $null
    "H81" = 2026.625
    "I81" = 2026.625
    "J81" = 0.0
    "K81" = 4053.25
    "L81" = 0.0
    "M81" = $null
    "N81" = -2992.25
    "H84" = 2026.625
    "I84" = 2026.625
    "J84" = 0.0
    "K84" = 20266.25
    "L84" = 0.0
    "M84" = $null
    "N84" = -14962.25
    "H97" = 40572.0
    "J97" = 40572.0
    "L97" = 40572.0
    "N97" = -42554.0
    "H100" = 891.9
    "I100" = 628.6667
    "K100" = 1257.3334
    "M100" = -716.3334
    "H103" = 0.0
    "J103" = 0.0
    "L103" = $null
    "N103" = 0.0
    "H122" = 68961.93
    "I122" = 2327.6667
    "J122" = 335499.0
    "K122" = 6983.000100000001
    "L122" = 1006497.0
    "M122" = -4533.000100000001
    "N122" = -1011397.0
    "H126" = 72421.484
    "I126" = 85849.56
    "J126" = 20947.166
    "K126" = 257548.68
    "L126" = 62841.49800000001
    "M126" = -255078.68
    "N126" = -67781.498
    "H132" = 33274.793
    "I132" = 34192.242
    "J132" = 2999.0
    "K132" = 102576.726
    "L132" = 8997.0
    "M132" = -100046.726
    "N132" = -14057.0
    "H136" = 3735.0
    "I136" = 3185.7407
    "K136" = 9557.222099999999
    "M136" = -7007.222099999999
}
foreach ($key in $changes.Keys) {
    $val = $changes[$key]
    if ($val -eq $null) {
        $ws.Range($key).ClearContents()
    } else {
        $ws.Range($key).Value = $val
    }
}
